# Qatar Stars League workbook update
# The source data rows got re-sorted/re-paired; for a number of adjacent
# row pairs the entire record (every column except the running index in
# column A) needs to be swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is a pair of adjacent worksheet rows whose B:AC content
# (id, HomeTeam, AwayTeam, score, odds, ...) must be exchanged, while
# column A (the sequential row index) stays where it is.
$rowPairs = @(
    @(18, 19),
    @(24, 25),
    @(29, 30),
    @(37, 38),
    @(45, 46),
    @(56, 57),
    @(62, 63),
    @(68, 69),
    @(75, 76),
    @(81, 82),
    @(83, 84),
    @(90, 91),
    @(94, 95),
    @(96, 97),
    @(98, 99),
    @(100, 101),
    @(108, 109),
    @(118, 119),
    @(120, 121)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B$r1" + ":AC$r1")
    $rangeB = $ws.Range("B$r2" + ":AC$r2")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
